# Insert a new row above row 472, shifting existing rows 472-532 down to 473-533,
# then populate the newly inserted row 472 with the new weekly data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(472).Insert()

$ws.Cells.Item(472, 1).Value = 3
$ws.Cells.Item(472, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(472, 3).Value = "Coquimbo"
$ws.Cells.Item(472, 4).Value = 45212
$ws.Cells.Item(472, 5).Value = 5
$ws.Cells.Item(472, 6).Value = 100112001
$ws.Cells.Item(472, 7).Value = "Berenjena"
$ws.Cells.Item(472, 8).Value = "Sin especificar"
$ws.Cells.Item(472, 9).Value = "Primera"
$ws.Cells.Item(472, 10).Value = 40
$ws.Cells.Item(472, 11).Value = 9000
$ws.Cells.Item(472, 12).Value = 9000
$ws.Cells.Item(472, 13).Value = 9000
$ws.Cells.Item(472, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(472, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(472, 16).Value = 150
$ws.Cells.Item(472, 17).Value = 60
$ws.Cells.Item(472, 18).Value = "Hortaliza"
